$d = $word.ActiveDocument

$replacements = @(
    @("164×9=", "322×4="),
    @("790×5=", "997×5="),
    @("905×4=", "590×9="),
    @("271×7=", "544×5="),
    @("596×8=", "529×9="),
    @("488×3=", "952×8="),
    @("354×8=", "403×7="),
    @("234×2=", "955×4="),
    @("933×5=", "102×7="),
    @("943×3=", "909×3="),
    @("317×7=", "606×6="),
    @("978×3=", "168×8="),
    @("511×8=", "654×3="),
    @("373×6=", "613×7="),
    @("290×2=", "826×5="),
    @("452×5=", "744×9="),
    @("304×3=", "452×3="),
    @("534×6=", "212×6="),
    @("794×8=", "108×7="),
    @("115×4=", "162×9="),
    @("924×3=", "723×4="),
    @("362×3=", "695×7="),
    @("366×6=", "806×6="),
    @("998×8=", "659×2="),
    @("234×5=", "530×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
